# Refresh the cryptos list data (price / 1h volume change columns, plus a
# couple of coin rows that swapped position) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.628.03'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.597.48'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.78%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '1.821.97'
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '1.567.47'
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.02'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.77'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').Value = '26.630.85'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '208.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.36%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.48'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('E28').Value = '  +0.65%  '
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0506'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E32').Value = '  -0.24%  '
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.623'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.28%  '
$ws.Range('D35').Value = '1.269.18'
$ws.Range('E35').Value = '  -1.85%  '
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.839'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('E40').Value = '  +2.64%  '
$ws.Range('E41').Value = '  +1.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.786'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '64.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.947'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +18.25%  '
$ws.Range('D45').Value = '1.735.04'
$ws.Range('E45').Value = '  +0.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.93'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('E48').Value = '  +4.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0508'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.61%  '
